$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header changes
$ws.Range("C1").Value = "rules"
$ws.Range("D1").Value = "adaptive_filter"

# Row 2
$ws.Range("D2").Value = "wRLS"
$ws.Range("E2").Value = 38.9818814560823
$ws.Range("F2").Value = 0.3062865873880715
$ws.Range("G2").Value = 29.97539273439184

# Row 3
$ws.Range("D3").Value = "wRLS"
$ws.Range("E3").Value = 37.33857140387859
$ws.Range("F3").Value = 0.2933748497009862
$ws.Range("G3").Value = 28.57544919503829

# Row 4
$ws.Range("D4").Value = "wRLS"
$ws.Range("E4").Value = 37.40029214320521
$ws.Range("F4").Value = 0.2938597989623685
$ws.Range("G4").Value = 28.61572899751141

# Row 5
$ws.Range("D5").Value = "wRLS"
$ws.Range("E5").Value = 39.23077416147777
$ws.Range("F5").Value = 0.308242175330825
$ws.Range("G5").Value = 30.195424658063

# Row 6
$ws.Range("D6").Value = "wRLS"
$ws.Range("E6").Value = 39.46205789280526
$ws.Range("F6").Value = 0.3100594068789351
$ws.Range("G6").Value = 30.33333559453169

# Row 7
$ws.Range("D7").Value = "wRLS"
$ws.Range("E7").Value = 39.33157713356496
$ws.Range("F7").Value = 0.3090341996550985
$ws.Range("G7").Value = 30.25194202452211

# Row 8
$ws.Range("D8").Value = "wRLS"
$ws.Range("E8").Value = 39.43717967675353
$ws.Range("F8").Value = 0.3098639349414563
$ws.Range("G8").Value = 30.35188308177086
